# Update experimental results (Power, Accuracy, Recall, F1 Score, Productivity)
# for the re-exported / reproduced run, per commit:
# "updating readme for reproducibility and re-exporting environment for
#  system cross-compatibility"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 312785.5599425451
$ws.Range("D2").Value = 0.2311661506707946
$ws.Range("E2").Value = 0.2244488977955912
$ws.Range("F2").Value = 0.2277580071174377
$ws.Range("G2").Value = 0.0007055872685981302

# Row 3
$ws.Range("C3").Value = 51613283.75073684
$ws.Range("D3").Value = 0.3906437425029988
$ws.Range("E3").Value = 0.3906958144494801
$ws.Range("F3").Value = 0.3906697767410863
$ws.Range("G3").Value = 0.00005679149090851102

# Row 4
$ws.Range("C4").Value = 25855353.23681396
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.3805651826179686
$ws.Range("F4").Value = 0.5513179492130926
$ws.Range("G4").Value = 0.00006087763453033134

# Row 5
$ws.Range("C5").Value = 51431834.82302178
$ws.Range("D5").Value = 0.7981553268279642
$ws.Range("E5").Value = 0.795921087709944
$ws.Range("F5").Value = 0.7970366415270639
$ws.Range("G5").Value = 0.000115932692966944

# Row 6
$ws.Range("C6").Value = 40845802.31418186
$ws.Range("D6").Value = 0.3906437425029988
$ws.Range("E6").Value = 0.3906958144494801
$ws.Range("F6").Value = 0.3906697767410863
$ws.Range("G6").Value = 0.00007176246196223316

# Row 7
$ws.Range("C7").Value = 25874295.03564626
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.7743268461743535
$ws.Range("F7").Value = 0.8728119600330554
$ws.Range("G7").Value = 0.0001959537320281384

# Row 8
$ws.Range("C8").Value = 40702246.41595799
$ws.Range("D8").Value = 0.7981553268279642
$ws.Range("E8").Value = 0.795921087709944
$ws.Range("F8").Value = 0.7970366415270639
$ws.Range("G8").Value = 0.0001464939073467507

# Row 9
$ws.Range("C9").Value = 20467060.82492377
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0.3805651826179686
$ws.Range("F9").Value = 0.5513179492130926
$ws.Range("G9").Value = 0.00007690467910695927

# Row 10
$ws.Range("C10").Value = 20482046.99586322
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.7743268461743535
$ws.Range("F10").Value = 0.8728119600330554
$ws.Range("G10").Value = 0.000247541892509867
